# Plantilla Cargue Contigencia-Cuba.xlsx
# Clears the sample/demo data rows (4-6) on the "Contingencia Entrada" sheet,
# leaving the cell formatting (styles) intact, resets the affected rows back
# to the default (auto) row height, and moves the active selection to F11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contingencia Entrada")

# Clear the sample data that was left in rows 4-6 (row 3 stays as the single
# example row). ClearContents keeps the cell styles/number formats applied
# to each cell, it just drops the stored values.
$ws.Range("A4:P6").ClearContents()

# The rows had an explicit custom height tied to the wrapped sample text;
# now that the text is gone, let Excel recompute (and drop) the custom
# row height so the rows fall back to the sheet's default height.
$ws.Range("A4:P6").EntireRow.AutoFit()

# Leave the selection where the editor left it.
$ws.Range("F11").Select()
